$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# New backlog item: cc clientes, fc en dolares pasa a la cuenta en pesos segun el tc
$ws.Range("A57").Value = "facturacion, cambiar de precio solo los articulos en dolares con t/c"
$ws.Range("B57").Value = "no comenzado"

# Update the view state to reflect scrolling down to the new row and the new selection
$win = $excel.ActiveWindow
$win.ScrollRow = 32
$win.ScrollColumn = 1
$ws.Range("A55").Select()
